# Auto-generated edit script: updates Leve profit-calculation columns (H-N)
# across multiple sheets per the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 478.2857
$ws.Range("I2").Value = 478.2857
$ws.Range("K2").Value = 478.2857
$ws.Range("M2").Value = -365.2857
$ws.Range("H19").Value = 551.2222
$ws.Range("I19").Value = 461.83334
$ws.Range("J19").Value = 730
$ws.Range("K19").Value = 461.83334
$ws.Range("L19").Value = 730
$ws.Range("M19").Value = -286.83334
$ws.Range("N19").Value = -1080
$ws.Range("H29").Value = 672.5454999999999
$ws.Range("J29").Value = 1666.6666
$ws.Range("L29").Value = 4999.9998
$ws.Range("N29").Value = -5561.9998
$ws.Range("H38").Value = 747.3333
$ws.Range("J38").Value = 2000
$ws.Range("L38").Value = 6000
$ws.Range("N38").Value = -6744
$ws.Range("H43").Value = 1847.75
$ws.Range("I43").Value = 1847.75
$ws.Range("K43").Value = 1847.75
$ws.Range("M43").Value = -1778.75
$ws.Range("H125").Value = 2366.7144
$ws.Range("I125").Value = 2199
$ws.Range("J125").Value = 2786
$ws.Range("K125").Value = 19791
$ws.Range("L125").Value = 25074
$ws.Range("M125").Value = -17331
$ws.Range("N125").Value = -29994
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2111.182
$ws.Range("J74").Value = 4000
$ws.Range("L74").Value = 4000
$ws.Range("N74").Value = -5748
$ws.Range("H77").Value = 2111.182
$ws.Range("J77").Value = 4000
$ws.Range("L77").Value = 20000
$ws.Range("N77").Value = -28736
$ws.Range("H97").Value = 5000
$ws.Range("I97").Value = 5000
$ws.Range("K97").Value = 5000
$ws.Range("M97").Value = -4504
$ws.Range("H110").Value = 5999.3335
$ws.Range("I110").Value = 3999
$ws.Range("J110").Value = 6999.5
$ws.Range("K110").Value = 3999
$ws.Range("L110").Value = 6999.5
$ws.Range("M110").Value = -1954
$ws.Range("N110").Value = -11089.5
$ws.Range("H122").Value = 10587.714
$ws.Range("I122").Value = 10138.529
$ws.Range("K122").Value = 30415.587
$ws.Range("M122").Value = -27965.587
$ws.Range("H132").Value = 1465.2858
$ws.Range("I132").Value = 1465.2858
$ws.Range("K132").Value = 4395.857400000001
$ws.Range("M132").Value = -1865.857400000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("H134").Value = 2799.2
$ws.Range("I134").Value = 2999.5
$ws.Range("K134").Value = 8998.5
$ws.Range("M134").Value = -6463.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1806.1
$ws.Range("I23").Value = 1894
$ws.Range("J23").Value = 1747.5
$ws.Range("K23").Value = 5682
$ws.Range("L23").Value = 5242.5
$ws.Range("M23").Value = -5447
$ws.Range("N23").Value = -5712.5
$ws.Range("H26").Value = 85
$ws.Range("I26").Value = 85
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 255
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 33
$ws.Range("N26").ClearContents()
$ws.Range("H75").Value = 916.3333
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 916.3333
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 2748.9999
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -4744.9999
$ws.Range("H78").Value = 916.3333
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 916.3333
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 8246.9997
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -18230.9997
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 782.44446
$ws.Range("I97").Value = 513.75
$ws.Range("K97").Value = 513.75
$ws.Range("M97").Value = -17.75
$ws.Range("H107").Value = 678.0909
$ws.Range("I107").Value = 584.44446
$ws.Range("K107").Value = 584.44446
$ws.Range("M107").Value = 1335.55554
$ws.Range("H132").Value = 2709
$ws.Range("I132").Value = 2891.4375
$ws.Range("K132").Value = 8674.3125
$ws.Range("M132").Value = -6144.3125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6545.727
$ws.Range("I7").Value = 6286.2856
$ws.Range("J7").Value = 6999.75
$ws.Range("K7").Value = 6286.2856
$ws.Range("L7").Value = 6999.75
$ws.Range("M7").Value = -6174.2856
$ws.Range("N7").Value = -7223.75
$ws.Range("H55").Value = 967.4
$ws.Range("I55").Value = 368.44446
$ws.Range("K55").Value = 368.44446
$ws.Range("M55").Value = -195.44446
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H100").Value = 5000
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H126").Value = 6545.727
$ws.Range("I126").Value = 6286.2856
$ws.Range("J126").Value = 6999.75
$ws.Range("K126").Value = 18858.8568
$ws.Range("L126").Value = 20999.25
$ws.Range("M126").Value = -16388.8568
$ws.Range("N126").Value = -25939.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2109.6
$ws.Range("I96").Value = 1934.3334
$ws.Range("J96").Value = 2372.5
$ws.Range("K96").Value = 1934.3334
$ws.Range("L96").Value = 2372.5
$ws.Range("M96").Value = -561.3334
$ws.Range("N96").Value = -5118.5
$ws.Range("H126").Value = 2346.6667
$ws.Range("I126").Value = 1906.6666
$ws.Range("J126").Value = 3666.6667
$ws.Range("K126").Value = 5719.9998
$ws.Range("L126").Value = 11000.0001
$ws.Range("M126").Value = -3249.9998
$ws.Range("N126").Value = -15940.0001
